$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 6.664768333333334
$ws.Cells.Item(2, 8).Value = 19.994305
$ws.Cells.Item(2, 9).Value = 0.06516174319532789
$ws.Cells.Item(2, 10).Value = 0.0651617431953279
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.237200333333333
$ws.Cells.Item(2, 14).Value = 6.711601
$ws.Cells.Item(2, 15).Value = 0.1121050933480713
$ws.Cells.Item(2, 16).Value = 0.1121050933480713
$ws.Cells.Item(2, 17).Value = 14.91042193692278
$ws.Cells.Item(2, 18).Value = 134.193797432305
$ws.Cells.Item(2, 19).Value = 0.007304963303635284
$ws.Cells.Item(2, 20).Value = 0.007304963303635286

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 6.664768333333334
$ws.Cells.Item(3, 8).Value = 19.994305
$ws.Cells.Item(3, 9).Value = 0.06516174319532789
$ws.Cells.Item(3, 10).Value = 0.0651617431953279
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.153682
$ws.Cells.Item(3, 14).Value = 9.461046
$ws.Cells.Item(3, 15).Value = 0.158029573718759
$ws.Cells.Item(3, 16).Value = 0.158029573718759
$ws.Cells.Item(3, 17).Value = 21.01855992700333
$ws.Cells.Item(3, 18).Value = 189.16703934303
$ws.Cells.Item(3, 19).Value = 0.01029748249992891
$ws.Cells.Item(3, 20).Value = 0.01029748249992891

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 6.664768333333334
$ws.Cells.Item(4, 8).Value = 19.994305
$ws.Cells.Item(4, 9).Value = 0.06516174319532789
$ws.Cells.Item(4, 10).Value = 0.0651617431953279
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 13.44189533333333
$ws.Cells.Item(4, 14).Value = 40.325686
$ws.Cells.Item(4, 15).Value = 0.6735672745377762
$ws.Cells.Item(4, 16).Value = 0.6735672745377762
$ws.Cells.Item(4, 17).Value = 89.58711835758113
$ws.Cells.Item(4, 18).Value = 806.2840652182301
$ws.Cells.Item(4, 19).Value = 0.04389081776820749
$ws.Cells.Item(4, 20).Value = 0.0438908177682075

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 6.664768333333334
$ws.Cells.Item(5, 8).Value = 19.994305
$ws.Cells.Item(5, 9).Value = 0.06516174319532789
$ws.Cells.Item(5, 10).Value = 0.0651617431953279
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.123499666666667
$ws.Cells.Item(5, 14).Value = 3.370499
$ws.Cells.Item(5, 15).Value = 0.05629805839539345
$ws.Cells.Item(5, 16).Value = 0.05629805839539345
$ws.Cells.Item(5, 17).Value = 7.487865000910557
$ws.Cells.Item(5, 18).Value = 67.39078500819501
$ws.Cells.Item(5, 19).Value = 0.003668479623556201
$ws.Cells.Item(5, 20).Value = 0.003668479623556202

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 47.25592399999999
$ws.Cells.Item(6, 8).Value = 141.767772
$ws.Cells.Item(6, 9).Value = 0.4620233187619072
$ws.Cells.Item(6, 10).Value = 0.4620233187619072
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 2.237200333333333
$ws.Cells.Item(6, 14).Value = 6.711601
$ws.Cells.Item(6, 15).Value = 0.1121050933480713
$ws.Cells.Item(6, 16).Value = 0.1121050933480713
$ws.Cells.Item(6, 17).Value = 105.7209689247747
$ws.Cells.Item(6, 18).Value = 951.4887203229719
$ws.Cells.Item(6, 19).Value = 0.05179516727878931
$ws.Cells.Item(6, 20).Value = 0.05179516727878932

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 47.25592399999999
$ws.Cells.Item(7, 8).Value = 141.767772
$ws.Cells.Item(7, 9).Value = 0.4620233187619072
$ws.Cells.Item(7, 10).Value = 0.4620233187619072
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.153682
$ws.Cells.Item(7, 14).Value = 9.461046
$ws.Cells.Item(7, 15).Value = 0.158029573718759
$ws.Cells.Item(7, 16).Value = 0.158029573718759
$ws.Cells.Item(7, 17).Value = 149.030156912168
$ws.Cells.Item(7, 18).Value = 1341.271412209512
$ws.Cells.Item(7, 19).Value = 0.0730133481120705
$ws.Cells.Item(7, 20).Value = 0.07301334811207053

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 47.25592399999999
$ws.Cells.Item(8, 8).Value = 141.767772
$ws.Cells.Item(8, 9).Value = 0.4620233187619072
$ws.Cells.Item(8, 10).Value = 0.4620233187619072
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 13.44189533333333
$ws.Cells.Item(8, 14).Value = 40.325686
$ws.Cells.Item(8, 15).Value = 0.6735672745377762
$ws.Cells.Item(8, 16).Value = 0.6735672745377762
$ws.Cells.Item(8, 17).Value = 635.2091842879546
$ws.Cells.Item(8, 18).Value = 5716.882658591592
$ws.Cells.Item(8, 19).Value = 0.311203787591356
$ws.Cells.Item(8, 20).Value = 0.311203787591356

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 47.25592399999999
$ws.Cells.Item(9, 8).Value = 141.767772
$ws.Cells.Item(9, 9).Value = 0.4620233187619072
$ws.Cells.Item(9, 10).Value = 0.4620233187619072
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.123499666666667
$ws.Cells.Item(9, 14).Value = 3.370499
$ws.Cells.Item(9, 15).Value = 0.05629805839539345
$ws.Cells.Item(9, 16).Value = 0.05629805839539345
$ws.Cells.Item(9, 17).Value = 53.09201486202533
$ws.Cells.Item(9, 18).Value = 477.828133758228
$ws.Cells.Item(9, 19).Value = 0.02601101577969133
$ws.Cells.Item(9, 20).Value = 0.02601101577969133

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 20.98736333333333
$ws.Cells.Item(10, 8).Value = 62.96209
$ws.Cells.Item(10, 9).Value = 0.2051944060881897
$ws.Cells.Item(10, 10).Value = 0.2051944060881898
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.237200333333333
$ws.Cells.Item(10, 14).Value = 6.711601
$ws.Cells.Item(10, 15).Value = 0.1121050933480713
$ws.Cells.Item(10, 16).Value = 0.1121050933480713
$ws.Cells.Item(10, 17).Value = 46.95293624512112
$ws.Cells.Item(10, 18).Value = 422.57642620609
$ws.Cells.Item(10, 19).Value = 0.02300333804901857
$ws.Cells.Item(10, 20).Value = 0.02300333804901857

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 20.98736333333333
$ws.Cells.Item(11, 8).Value = 62.96209
$ws.Cells.Item(11, 9).Value = 0.2051944060881897
$ws.Cells.Item(11, 10).Value = 0.2051944060881898
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.153682
$ws.Cells.Item(11, 14).Value = 9.461046
$ws.Cells.Item(11, 15).Value = 0.158029573718759
$ws.Cells.Item(11, 16).Value = 0.158029573718759
$ws.Cells.Item(11, 17).Value = 66.18746997179333
$ws.Cells.Item(11, 18).Value = 595.68722974614
$ws.Cells.Item(11, 19).Value = 0.03242678452359055
$ws.Cells.Item(11, 20).Value = 0.03242678452359056

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 20.98736333333333
$ws.Cells.Item(12, 8).Value = 62.96209
$ws.Cells.Item(12, 9).Value = 0.2051944060881897
$ws.Cells.Item(12, 10).Value = 0.2051944060881898
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 13.44189533333333
$ws.Cells.Item(12, 14).Value = 40.325686
$ws.Cells.Item(12, 15).Value = 0.6735672745377762
$ws.Cells.Item(12, 16).Value = 0.6735672745377762
$ws.Cells.Item(12, 17).Value = 282.1099412493045
$ws.Cells.Item(12, 18).Value = 2538.98947124374
$ws.Cells.Item(12, 19).Value = 0.1382122368592196
$ws.Cells.Item(12, 20).Value = 0.1382122368592197

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 20.98736333333333
$ws.Cells.Item(13, 8).Value = 62.96209
$ws.Cells.Item(13, 9).Value = 0.2051944060881897
$ws.Cells.Item(13, 10).Value = 0.2051944060881898
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.123499666666667
$ws.Cells.Item(13, 14).Value = 3.370499
$ws.Cells.Item(13, 15).Value = 0.05629805839539345
$ws.Cells.Item(13, 16).Value = 0.05629805839539345
$ws.Cells.Item(13, 17).Value = 23.57929570921223
$ws.Cells.Item(13, 18).Value = 212.21366138291
$ws.Cells.Item(13, 19).Value = 0.01155204665636098
$ws.Cells.Item(13, 20).Value = 0.01155204665636098

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 27.37233166666667
$ws.Cells.Item(14, 8).Value = 82.116995
$ws.Cells.Item(14, 9).Value = 0.2676205319545753
$ws.Cells.Item(14, 10).Value = 0.2676205319545753
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 2.237200333333333
$ws.Cells.Item(14, 14).Value = 6.711601
$ws.Cells.Item(14, 15).Value = 0.1121050933480713
$ws.Cells.Item(14, 16).Value = 0.1121050933480713
$ws.Cells.Item(14, 17).Value = 61.23738952877723
$ws.Cells.Item(14, 18).Value = 551.136505758995
$ws.Cells.Item(14, 19).Value = 0.03000162471662817
$ws.Cells.Item(14, 20).Value = 0.03000162471662817

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 27.37233166666667
$ws.Cells.Item(15, 8).Value = 82.116995
$ws.Cells.Item(15, 9).Value = 0.2676205319545753
$ws.Cells.Item(15, 10).Value = 0.2676205319545753
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 3.153682
$ws.Cells.Item(15, 14).Value = 9.461046
$ws.Cells.Item(15, 15).Value = 0.158029573718759
$ws.Cells.Item(15, 16).Value = 0.158029573718759
$ws.Cells.Item(15, 17).Value = 86.32362967519667
$ws.Cells.Item(15, 18).Value = 776.91266707677
$ws.Cells.Item(15, 19).Value = 0.04229195858316906
$ws.Cells.Item(15, 20).Value = 0.04229195858316907

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 27.37233166666667
$ws.Cells.Item(16, 8).Value = 82.116995
$ws.Cells.Item(16, 9).Value = 0.2676205319545753
$ws.Cells.Item(16, 10).Value = 0.2676205319545753
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 13.44189533333333
$ws.Cells.Item(16, 14).Value = 40.325686
$ws.Cells.Item(16, 15).Value = 0.6735672745377762
$ws.Cells.Item(16, 16).Value = 0.6735672745377762
$ws.Cells.Item(16, 17).Value = 367.9360172926189
$ws.Cells.Item(16, 18).Value = 3311.424155633571
$ws.Cells.Item(16, 19).Value = 0.1802604323189931
$ws.Cells.Item(16, 20).Value = 0.1802604323189931

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 27.37233166666667
$ws.Cells.Item(17, 8).Value = 82.116995
$ws.Cells.Item(17, 9).Value = 0.2676205319545753
$ws.Cells.Item(17, 10).Value = 0.2676205319545753
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.123499666666667
$ws.Cells.Item(17, 14).Value = 3.370499
$ws.Cells.Item(17, 15).Value = 0.05629805839539345
$ws.Cells.Item(17, 16).Value = 0.05629805839539345
$ws.Cells.Item(17, 17).Value = 30.75280550338945
$ws.Cells.Item(17, 18).Value = 276.775249530505
$ws.Cells.Item(17, 19).Value = 0.01506651633578494
$ws.Cells.Item(17, 20).Value = 0.01506651633578494
